$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction9")

# Clear out columns C through R in row 1 (they will no longer hold data)
$ws.Range("C1:R1").ClearContents()

# Update A1/B1 with the new values
$ws.Range("A1").Value = 18
$ws.Range("B1").Value = 19
